$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (the second game/event record) entirely; row 2's former
# neighbour is gone and row indices shift up.
$ws.Rows.Item(3).Delete()

# Update remaining row 2 with the new record's data.
# id - stored as text "0" (not a number)
$ws.Range("A2").Value = "'0"

# host
$ws.Range("B2").Value = "Nic Bolton"

# date - stored as text "20230906" (not a number)
$ws.Range("C2").Value = "'20230906"

# time - stored as text "1000" (not a number)
$ws.Range("D2").Value = "'1000"

# location
$ws.Range("E2").Value = "Prospects Athletics"

# player limit (numeric)
$ws.Range("F2").Value = 5

# goalie limit column is no longer populated for this record
$ws.Range("G2").ClearContents()

# goalie limit (numeric) stays the same
$ws.Range("H2").Value = 2

# goalie list
$ws.Range("I2").Value = "0;"

# Strip the quote-prefix formatting picked up from the text-forcing
# assignments above so the cells keep the workbook's default style.
$ws.Range("A2:D2").Style = "Normal"
